# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font, border, centered alignment) from an
# existing header cell onto the three new header cells so they match the
# rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Every player row shares the team's season record.
for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 29).Value = 55
    $ws.Cells.Item($r, 30).Value = 58
    $ws.Cells.Item($r, 31).Value = 0
}
